# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.482.10"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.894.99"
$ws.Range("E3").Value = "  -1.18%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.22"
$ws.Range("E5").Value = "  +3.74%  "

$ws.Range("E6").Value = "  -0.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4759"
$ws.Range("E7").Value = "  -1.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3998"
$ws.Range("E8").Value = "  -1.51%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.28"
$ws.Range("E9").Value = "  -1.38%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08018"
$ws.Range("E10").Value = "  -2.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9898"
$ws.Range("E11").Value = "  -1.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.14"
$ws.Range("E12").Value = "  -0.92%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.915.10"
$ws.Range("E13").Value = "  -0.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.930"
$ws.Range("E14").Value = "  -2.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.079"
$ws.Range("E15").Value = "  -2.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.05"
$ws.Range("E16").Value = "  -2.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06806"
$ws.Range("E17").Value = "  -0.98%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("E18").Value = "  -0.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001019"
$ws.Range("E19").Value = "  -1.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.31"
$ws.Range("E20").Value = "  -1.63%  "

$ws.Range("E21").Value = "  -0.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.481.50"
$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.508"
$ws.Range("E23").Value = "  -2.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.65"
$ws.Range("E24").Value = "  -1.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.151"
$ws.Range("E25").Value = "  -1.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.149.86"
$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.30"
$ws.Range("E27").Value = "  +0.96%  "

$ws.Range("E28").Value = "  -1.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.64"
$ws.Range("E29").Value = "  -1.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.050"
$ws.Range("E30").Value = "  -2.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.95"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9958"
$ws.Range("E32").Value = "  -1.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09550"
$ws.Range("E33").Value = "  -0.91%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.472"
$ws.Range("E34").Value = "  -2.69%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.528"
$ws.Range("E35").Value = "  -0.69%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.382"
$ws.Range("E36").Value = "  +0.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06424"
$ws.Range("E37").Value = "  +5.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02240"
$ws.Range("E38").Value = "  -1.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.197"
$ws.Range("E39").Value = "  +1.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5829"
$ws.Range("E40").Value = "  -2.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.54"
$ws.Range("E41").Value = "  -3.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.743"
$ws.Range("E42").Value = "  -3.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1820"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.426"
$ws.Range("E44").Value = "  +2.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.230"
$ws.Range("E45").Value = "  -3.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.18"
$ws.Range("E46").Value = "  -1.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5494"
$ws.Range("E47").Value = "  -1.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.07332"
$ws.Range("E48").Value = "  -3.64%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.951"
$ws.Range("E49").Value = "  +0.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "116.37"
$ws.Range("E50").Value = "  -2.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.370"
$ws.Range("E51").Value = "  -2.20%  "
